$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 118
$ws1.Range("F5").Value = 5402
$ws1.Range("F6").Value = 73
$ws1.Range("F7").Value = 890
$ws1.Range("F8").Value = 134
$ws1.Range("F9").Value = 2397
$ws1.Range("F11").Value = 54
$ws1.Range("F12").Value = 2249
$ws1.Range("F13").Value = 57

# Sheet 4: 全部类型 (All Types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 118
$ws4.Range("F5").Value = 5402
$ws4.Range("F7").Value = 73
$ws4.Range("F9").Value = 890
$ws4.Range("F10").Value = 134
$ws4.Range("F11").Value = 2397
$ws4.Range("F14").Value = 54
$ws4.Range("F15").Value = 2249
$ws4.Range("F16").Value = 57
